$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as Text so values like "1.015"
# or "28.306.81" are not auto-converted into numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = '28.306.81'
$ws.Range("E2").Value2 = '  +0.72%  '

$ws.Range("D3").Value2 = '1.864.39'
$ws.Range("E3").Value2 = '  -0.47%  '

$ws.Range("E4").Value2 = '  +1.46%  '

$ws.Range("D5").Value2 = '315.81'
$ws.Range("E5").Value2 = '  +0.71%  '

$ws.Range("D6").Value2 = '1.015'
$ws.Range("E6").Value2 = '  +1.22%  '

$ws.Range("D7").Value2 = '0.5093'
$ws.Range("E7").Value2 = '  -0.28%  '

$ws.Range("D8").Value2 = '0.3940'
$ws.Range("E8").Value2 = '  +1.14%  '

$ws.Range("D9").Value2 = '0.08440'
$ws.Range("E9").Value2 = '  +1.15%  '

$ws.Range("D10").Value2 = '1.106'
$ws.Range("E10").Value2 = '  -1.03%  '

$ws.Range("D11").Value2 = '6.227'
$ws.Range("E11").Value2 = '  +0.04%  '

$ws.Range("D12").Value2 = '20.40'
$ws.Range("E12").Value2 = '  -0.57%  '

$ws.Range("D13").Value2 = '1.812.57'
$ws.Range("E13").Value2 = '  -3.18%  '

$ws.Range("B14").Value2 = 'Chainlink'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value2 = '7.202'
$ws.Range("E14").Value2 = '  -0.47%  '

$ws.Range("B15").Value2 = 'BinanceUSD'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value2 = '1.016'
$ws.Range("E15").Value2 = '  +1.63%  '

$ws.Range("D16").Value2 = '0.00001108'
$ws.Range("E16").Value2 = '  +0.84%  '

$ws.Range("D17").Value2 = '90.40'
$ws.Range("E17").Value2 = '  -0.47%  '

$ws.Range("D18").Value2 = '0.06753'
$ws.Range("E18").Value2 = '  +1.09%  '

$ws.Range("D19").Value2 = '17.64'
$ws.Range("E19").Value2 = '  -0.51%  '

$ws.Range("D20").Value2 = '1.014'
$ws.Range("E20").Value2 = '  +1.20%  '

$ws.Range("D21").Value2 = '5.928'
$ws.Range("E21").Value2 = '  -1.31%  '

$ws.Range("D22").Value2 = '28.324.31'
$ws.Range("E22").Value2 = '  +0.59%  '

$ws.Range("D23").Value2 = '11.11'
$ws.Range("E23").Value2 = '  +0.10%  '

$ws.Range("D24").Value2 = '2.278'
$ws.Range("E24").Value2 = '  +0.65%  '

$ws.Range("D25").Value2 = '161.82'
$ws.Range("E25").Value2 = '  +1.89%  '

$ws.Range("D26").Value2 = '2.036.63'
$ws.Range("E26").Value2 = '  -2.71%  '

$ws.Range("D27").Value2 = '20.72'
$ws.Range("E27").Value2 = '  +0.17%  '

$ws.Range("D28").Value2 = '2.352'
$ws.Range("E28").Value2 = '  -3.92%  '

$ws.Range("D29").Value2 = '125.89'
$ws.Range("E29").Value2 = '  -0.18%  '

$ws.Range("D30").Value2 = '0.1047'
$ws.Range("E30").Value2 = '  -0.56%  '

$ws.Range("D31").Value2 = '1.033'
$ws.Range("E31").Value2 = '  -0.04%  '

$ws.Range("D32").Value2 = '5.749'
$ws.Range("E32").Value2 = '  -1.67%  '

$ws.Range("D33").Value2 = '3.628'
$ws.Range("E33").Value2 = '  +0.65%  '

$ws.Range("D34").Value2 = '0.02424'
$ws.Range("E34").Value2 = '  -0.60%  '

$ws.Range("D35").Value2 = '0.06436'
$ws.Range("E35").Value2 = '  -1.67%  '

$ws.Range("D36").Value2 = '0.2179'
$ws.Range("E36").Value2 = '  -1.34%  '

$ws.Range("D37").Value2 = '8.803'
$ws.Range("E37").Value2 = '  -7.45%  '

$ws.Range("D38").Value2 = '1.260'
$ws.Range("E38").Value2 = '  +1.41%  '

$ws.Range("D39").Value2 = '1.176'
$ws.Range("E39").Value2 = '  -1.28%  '

$ws.Range("D40").Value2 = '0.6356'
$ws.Range("E40").Value2 = '  -1.53%  '

$ws.Range("D41").Value2 = '4.967'
$ws.Range("E41").Value2 = '  -0.21%  '

$ws.Range("D42").Value2 = '11.21'
$ws.Range("E42").Value2 = '  -0.66%  '

$ws.Range("D43").Value2 = '0.5998'
$ws.Range("E43").Value2 = '  -1.18%  '

$ws.Range("D44").Value2 = '12.98'
$ws.Range("E44").Value2 = '  -0.39%  '

$ws.Range("D45").Value2 = '3.691'
$ws.Range("E45").Value2 = '  +0.34%  '

$ws.Range("D46").Value2 = '1.213'
$ws.Range("E46").Value2 = '  -4.71%  '

$ws.Range("D47").Value2 = '1.980'
$ws.Range("E47").Value2 = '  -1.30%  '

$ws.Range("D48").Value2 = '1.198'
$ws.Range("E48").Value2 = '  -2.54%  '

$ws.Range("D49").Value2 = '120.79'
$ws.Range("E49").Value2 = '  +0.21%  '

$ws.Range("D50").Value2 = '0.06822'
$ws.Range("E50").Value2 = '  -1.06%  '

$ws.Range("D51").Value2 = '76.05'
$ws.Range("E51").Value2 = '  -2.11%  '

# Restore the default "Normal" style on column D so no stray number-format
# style is left applied to the cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
